# Mark the "版本信息" (version info) TODO item as done and set its
# resolution date, then move the active selection to D25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C24: status -> 完成 checkmark "√" (matches the other "done" rows)
$ws.Range("C24").Value = "√"

# D24: resolution date (2017-03-15 -> serial 42809)
$ws.Range("D24").Value = 42809

# Move the selection to D25, matching the updated view state
$ws.Range("D25").Select()
